# Re-shuffle the "observation" columns (A,B,E,F,G,H,K:N,Q,R,AC) between rows
# 2..13 according to a fixed permutation, while leaving the row-local columns
# (C,D,I,P,S,T,U,V,W,Y,Z,AA,AB,AD,AE,AG,AT,AW,AX,AY, ...) untouched.
#
# Mapping: new row R gets the "moving" data that currently sits in row Map[R].

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$map = @{
    2  = 13
    3  = 2
    4  = 8
    5  = 12
    6  = 3
    7  = 9
    8  = 6
    9  = 5
    10 = 11
    11 = 4
    12 = 7
    13 = 10
}

$firstRow = 2
$lastRow = 13

# Snapshot the values that move, for every row, before anything is written
# (so that writes to one row never clobber data still needed for another).
$snapA  = @{}
$snapB  = @{}
$snapE  = @{}
$snapF  = @{}
$snapG  = @{}
$snapH  = @{}
$snapQ  = @{}
$snapR  = @{}
$snapAC = @{}
$snapHasKLMN = @{}

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $snapA[$r]  = $ws.Range("A$r").Value2
    $snapB[$r]  = $ws.Range("B$r").Value2
    $snapE[$r]  = $ws.Range("E$r").Value2
    $snapF[$r]  = $ws.Range("F$r").Value2
    $snapG[$r]  = $ws.Range("G$r").Value2
    $snapH[$r]  = $ws.Range("H$r").Value2
    $snapQ[$r]  = $ws.Range("Q$r").Value2
    $snapR[$r]  = $ws.Range("R$r").Value2
    $snapAC[$r] = $ws.Range("AC$r").Value2
    # K:N (Ålder-Stadium, Kön, Aktivitet, Metod) are present (as blank cells)
    # exactly for the "Tretåig hackspett" (TaxonId 100109 / B = 56398) rows.
    $snapHasKLMN[$r] = ($snapB[$r] -eq 56398)
}

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $src = $map[$r]

    $ws.Range("A$r").Value = $snapA[$src]
    $ws.Range("B$r").Value = $snapB[$src]
    $ws.Range("E$r").Value = $snapE[$src]
    $ws.Range("F$r").Value = $snapF[$src]
    $ws.Range("G$r").Value = $snapG[$src]
    $ws.Range("H$r").Value = $snapH[$src]
    $ws.Range("Q$r").Value = $snapQ[$src]
    $ws.Range("R$r").Value = $snapR[$src]

    if ($snapAC[$src] -eq $null -or $snapAC[$src] -eq "") {
        $ws.Range("AC$r").ClearContents()
    } else {
        $ws.Range("AC$r").Value = $snapAC[$src]
    }

    if ($snapHasKLMN[$src]) {
        $ws.Range("K$r").Value = ""
        $ws.Range("L$r").Value = ""
        $ws.Range("M$r").Value = ""
        $ws.Range("N$r").Value = ""
    } else {
        $ws.Range("K$r").ClearContents()
        $ws.Range("L$r").ClearContents()
        $ws.Range("M$r").ClearContents()
        $ws.Range("N$r").ClearContents()
    }
}
